$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date value for all data rows: 45184 -> 45186
for ($i = 2; $i -le 431; $i++) {
    $ws.Cells.Item($i, 3).Value2 = 45186
}

# Add a display-text second argument to existing HYPERLINK formulas
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 34151-2019.xlsx", "A 34151-2019")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 34151-2019.png", "A 34151-2019")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/knärot/A 34151-2019.png", "A 34151-2019")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 34151-2019.docx", "A 34151-2019")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 34151-2019.docx", "A 34151-2019")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 34151-2019.docx", "A 34151-2019")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 34151-2019.docx", "A 34151-2019")'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 37759-2020.xlsx", "A 37759-2020")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 37759-2020.png", "A 37759-2020")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 37759-2020.docx", "A 37759-2020")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 37759-2020.docx", "A 37759-2020")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 37759-2020.docx", "A 37759-2020")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 37759-2020.docx", "A 37759-2020")'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 34968-2020.xlsx", "A 34968-2020")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 34968-2020.png", "A 34968-2020")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 34968-2020.docx", "A 34968-2020")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 34968-2020.docx", "A 34968-2020")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 34968-2020.docx", "A 34968-2020")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 34968-2020.docx", "A 34968-2020")'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 11177-2020.xlsx", "A 11177-2020")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 11177-2020.png", "A 11177-2020")'
$ws.Range("U5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/knärot/A 11177-2020.png", "A 11177-2020")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 11177-2020.docx", "A 11177-2020")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 11177-2020.docx", "A 11177-2020")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 11177-2020.docx", "A 11177-2020")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 11177-2020.docx", "A 11177-2020")'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 13012-2022.xlsx", "A 13012-2022")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 13012-2022.png", "A 13012-2022")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 13012-2022.docx", "A 13012-2022")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 13012-2022.docx", "A 13012-2022")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 13012-2022.docx", "A 13012-2022")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 13012-2022.docx", "A 13012-2022")'
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 34145-2019.xlsx", "A 34145-2019")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 34145-2019.png", "A 34145-2019")'
$ws.Range("U7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/knärot/A 34145-2019.png", "A 34145-2019")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 34145-2019.docx", "A 34145-2019")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 34145-2019.docx", "A 34145-2019")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 34145-2019.docx", "A 34145-2019")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 34145-2019.docx", "A 34145-2019")'
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 37380-2020.xlsx", "A 37380-2020")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 37380-2020.png", "A 37380-2020")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 37380-2020.docx", "A 37380-2020")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 37380-2020.docx", "A 37380-2020")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 37380-2020.docx", "A 37380-2020")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 37380-2020.docx", "A 37380-2020")'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 66048-2020.xlsx", "A 66048-2020")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 66048-2020.png", "A 66048-2020")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 66048-2020.docx", "A 66048-2020")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 66048-2020.docx", "A 66048-2020")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 66048-2020.docx", "A 66048-2020")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 66048-2020.docx", "A 66048-2020")'
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 21966-2022.xlsx", "A 21966-2022")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 21966-2022.png", "A 21966-2022")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 21966-2022.docx", "A 21966-2022")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 21966-2022.docx", "A 21966-2022")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 21966-2022.docx", "A 21966-2022")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 21966-2022.docx", "A 21966-2022")'
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 46853-2018.xlsx", "A 46853-2018")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 46853-2018.png", "A 46853-2018")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 46853-2018.docx", "A 46853-2018")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 46853-2018.docx", "A 46853-2018")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 46853-2018.docx", "A 46853-2018")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 46853-2018.docx", "A 46853-2018")'
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 46863-2018.xlsx", "A 46863-2018")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 46863-2018.png", "A 46863-2018")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 46863-2018.docx", "A 46863-2018")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 46863-2018.docx", "A 46863-2018")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 46863-2018.docx", "A 46863-2018")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 46863-2018.docx", "A 46863-2018")'
$ws.Range("S13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 12671-2022.xlsx", "A 12671-2022")'
$ws.Range("T13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 12671-2022.png", "A 12671-2022")'
$ws.Range("V13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 12671-2022.docx", "A 12671-2022")'
$ws.Range("W13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 12671-2022.docx", "A 12671-2022")'
$ws.Range("X13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 12671-2022.docx", "A 12671-2022")'
$ws.Range("Y13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 12671-2022.docx", "A 12671-2022")'
$ws.Range("S14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 44717-2022.xlsx", "A 44717-2022")'
$ws.Range("T14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 44717-2022.png", "A 44717-2022")'
$ws.Range("V14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 44717-2022.docx", "A 44717-2022")'
$ws.Range("W14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 44717-2022.docx", "A 44717-2022")'
$ws.Range("X14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 44717-2022.docx", "A 44717-2022")'
$ws.Range("Y14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 44717-2022.docx", "A 44717-2022")'
$ws.Range("S15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 8638-2019.xlsx", "A 8638-2019")'
$ws.Range("T15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 8638-2019.png", "A 8638-2019")'
$ws.Range("V15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 8638-2019.docx", "A 8638-2019")'
$ws.Range("W15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 8638-2019.docx", "A 8638-2019")'
$ws.Range("X15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 8638-2019.docx", "A 8638-2019")'
$ws.Range("Y15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 8638-2019.docx", "A 8638-2019")'
$ws.Range("S16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 20041-2020.xlsx", "A 20041-2020")'
$ws.Range("T16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 20041-2020.png", "A 20041-2020")'
$ws.Range("V16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 20041-2020.docx", "A 20041-2020")'
$ws.Range("W16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 20041-2020.docx", "A 20041-2020")'
$ws.Range("X16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 20041-2020.docx", "A 20041-2020")'
$ws.Range("Y16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 20041-2020.docx", "A 20041-2020")'
$ws.Range("S17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 23515-2020.xlsx", "A 23515-2020")'
$ws.Range("T17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 23515-2020.png", "A 23515-2020")'
$ws.Range("V17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 23515-2020.docx", "A 23515-2020")'
$ws.Range("W17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 23515-2020.docx", "A 23515-2020")'
$ws.Range("X17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 23515-2020.docx", "A 23515-2020")'
$ws.Range("Y17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 23515-2020.docx", "A 23515-2020")'
$ws.Range("S18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 58011-2020.xlsx", "A 58011-2020")'
$ws.Range("T18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 58011-2020.png", "A 58011-2020")'
$ws.Range("V18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 58011-2020.docx", "A 58011-2020")'
$ws.Range("W18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 58011-2020.docx", "A 58011-2020")'
$ws.Range("X18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 58011-2020.docx", "A 58011-2020")'
$ws.Range("Y18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 58011-2020.docx", "A 58011-2020")'
$ws.Range("S19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 64960-2020.xlsx", "A 64960-2020")'
$ws.Range("T19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 64960-2020.png", "A 64960-2020")'
$ws.Range("U19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/knärot/A 64960-2020.png", "A 64960-2020")'
$ws.Range("V19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 64960-2020.docx", "A 64960-2020")'
$ws.Range("W19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 64960-2020.docx", "A 64960-2020")'
$ws.Range("X19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 64960-2020.docx", "A 64960-2020")'
$ws.Range("Y19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 64960-2020.docx", "A 64960-2020")'
$ws.Range("S20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 63301-2021.xlsx", "A 63301-2021")'
$ws.Range("T20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 63301-2021.png", "A 63301-2021")'
$ws.Range("V20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 63301-2021.docx", "A 63301-2021")'
$ws.Range("W20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 63301-2021.docx", "A 63301-2021")'
$ws.Range("X20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 63301-2021.docx", "A 63301-2021")'
$ws.Range("Y20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 63301-2021.docx", "A 63301-2021")'
$ws.Range("S21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 65330-2021.xlsx", "A 65330-2021")'
$ws.Range("T21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 65330-2021.png", "A 65330-2021")'
$ws.Range("V21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 65330-2021.docx", "A 65330-2021")'
$ws.Range("W21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 65330-2021.docx", "A 65330-2021")'
$ws.Range("X21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 65330-2021.docx", "A 65330-2021")'
$ws.Range("Y21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 65330-2021.docx", "A 65330-2021")'
$ws.Range("S22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 24454-2022.xlsx", "A 24454-2022")'
$ws.Range("T22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 24454-2022.png", "A 24454-2022")'
$ws.Range("V22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 24454-2022.docx", "A 24454-2022")'
$ws.Range("W22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 24454-2022.docx", "A 24454-2022")'
$ws.Range("X22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 24454-2022.docx", "A 24454-2022")'
$ws.Range("Y22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 24454-2022.docx", "A 24454-2022")'
$ws.Range("S23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 33751-2022.xlsx", "A 33751-2022")'
$ws.Range("T23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 33751-2022.png", "A 33751-2022")'
$ws.Range("V23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 33751-2022.docx", "A 33751-2022")'
$ws.Range("W23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 33751-2022.docx", "A 33751-2022")'
$ws.Range("X23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 33751-2022.docx", "A 33751-2022")'
$ws.Range("Y23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 33751-2022.docx", "A 33751-2022")'
$ws.Range("S24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 41511-2022.xlsx", "A 41511-2022")'
$ws.Range("T24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 41511-2022.png", "A 41511-2022")'
$ws.Range("V24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 41511-2022.docx", "A 41511-2022")'
$ws.Range("W24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 41511-2022.docx", "A 41511-2022")'
$ws.Range("X24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 41511-2022.docx", "A 41511-2022")'
$ws.Range("Y24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 41511-2022.docx", "A 41511-2022")'
$ws.Range("S25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 692-2023.xlsx", "A 692-2023")'
$ws.Range("T25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 692-2023.png", "A 692-2023")'
$ws.Range("V25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 692-2023.docx", "A 692-2023")'
$ws.Range("W25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 692-2023.docx", "A 692-2023")'
$ws.Range("X25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 692-2023.docx", "A 692-2023")'
$ws.Range("Y25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 692-2023.docx", "A 692-2023")'
$ws.Range("S26").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/artfynd/A 2303-2023.xlsx", "A 2303-2023")'
$ws.Range("T26").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/kartor/A 2303-2023.png", "A 2303-2023")'
$ws.Range("V26").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 2303-2023.docx", "A 2303-2023")'
$ws.Range("W26").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 2303-2023.docx", "A 2303-2023")'
$ws.Range("X26").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 2303-2023.docx", "A 2303-2023")'
$ws.Range("Y26").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 2303-2023.docx", "A 2303-2023")'
$ws.Range("U175").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/knärot/A 57642-2020.png", "A 57642-2020")'
$ws.Range("V175").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomål/A 57642-2020.docx", "A 57642-2020")'
$ws.Range("W175").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/klagomålsmail/A 57642-2020.docx", "A 57642-2020")'
$ws.Range("X175").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsyn/A 57642-2020.docx", "A 57642-2020")'
$ws.Range("Y175").Formula = '=HYPERLINK("https://klasma.github.io/Logging_STRANGNAS/tillsynsmail/A 57642-2020.docx", "A 57642-2020")'
